# Xenium FOV import: rename the "brain1" sample to "brain" and drop the
# second ("brain2") sample row entirely, then leave the selection on the
# former brain2 row's position (now B3, after the row shift).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 held the "brain2" sample; remove it (all cells shift up one row).
$ws.Rows.Item(3).Delete()

# Row 2 ("brain1") becomes the lone sample; rename it to "brain".
$ws.Range("B2").Value = "brain"

# Update the active selection to B3 (matches the post-edit cursor position).
$ws.Range("B3").Select() | Out-Null
